# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values on Sheet1 for rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 1
